$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy cell formats from analogous existing rows first -----------------
$ws.Range("A9").Copy()
$ws.Range("A10:A11").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("C9:D9").Copy()
$ws.Range("C10:D10").PasteSpecial(-4122)
$ws.Range("C9:D9").Copy()
$ws.Range("C11:D11").PasteSpecial(-4122)

$ws.Range("E8").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("E8:F8").Copy()
$ws.Range("E11:F11").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Rows.Item(10).RowHeight = 15.75
$ws.Rows.Item(11).RowHeight = 15.75

# --- Row 10 : maternelle2 / Louise -----------------------------------------
$ws.Range("A10").Value2 = 43435.669681574072
$ws.Range("B10").Value2 = 43429
$ws.Range("B10").NumberFormat = "m/d/yyyy"
$ws.Range("C10").Value = "maternelle2"
$ws.Range("D10").Value = "Louise"
$ws.Range("E10").Value2 = 0.375
$ws.Range("I10").Value2 = 2

# --- Row 11 : maternelle2 / Joséphine ---------------------------------------
$ws.Range("A11").Value2 = 43435.669681574072
$ws.Range("B11").Value2 = 43429
$ws.Range("B11").NumberFormat = "m/d/yyyy"
$ws.Range("C11").Value = "maternelle2"
$ws.Range("D11").Value = "Joséphine"
$ws.Range("E11").Value2 = 0.375
$ws.Range("F11").Value2 = 0.72916666666424135
$ws.Range("H11").Value2 = 1

$ws.Range("I10").Select()
